# The workbook originally lists its two worksheets in the order:
#   1) "2022-Q2"  (quarterly fund-holding detail table)
#   2) "总计"      (summary / totals table)
#
# This edit re-sorts the sheet tabs so that the summary sheet "总计" comes
# first and the quarterly detail sheet "2022-Q2" comes second, while
# keeping "2022-Q2" as the active/selected tab (as it was originally).

$wb = $excel.ActiveWorkbook

# Move "总计" to be the first sheet in the workbook (before whatever sheet
# currently occupies position 1).
$totalSheet = $wb.Worksheets.Item("总计")
$beforeTarget = $wb.Worksheets.Item(1)
$totalSheet.Move($beforeTarget)

# Re-acquire a fresh reference to "2022-Q2" now that the sheet collection
# has been reordered, and make it the active tab again.
$quarterSheet = $wb.Worksheets.Item("2022-Q2")
$quarterSheet.Activate()
